# Atualização automática SALDO_PECAS (17/11/2025 15:53)
# Adds a new data row (row 12) to the PRINCIPAL sheet, mirroring the
# existing rows' layout (columns A-M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

$ws.Cells.Item($row, 1).Value()  = "DF"                           # A - UF
$ws.Cells.Item($row, 2).Value()  = "DF17110"                      # B - FRU
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value()  = "125200"                       # C - SUB1 (kept as text)
$ws.Cells.Item($row, 6).Value()  = "T"                             # F - DESCRICAO
$ws.Cells.Item($row, 7).Value()  = "T"                             # G - MAQUINAS
$ws.Cells.Item($row, 8).Value()  = "T - (T 30/11/26_24H) - DF"    # H - CLIENTE
$ws.Cells.Item($row, 9).Value()  = "30/11/26"                      # I - DATA_FIM
$ws.Cells.Item($row, 10).Value() = "24H"                           # J - SLA
$ws.Cells.Item($row, 11).Value() = "17/11/25"                      # K - DATA_VERIFICACAO
$ws.Cells.Item($row, 12).Value() = "DENTRO"                        # L - STATUS
